$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DeviceList")
$ws.Activate()

# Delete entire column D (shifts E:I left to D:H)
$ws.Columns("D").Delete()

# The conditional formatting on B2:I2 should shrink to B2:H2 now that
# the range lost one column.
for ($i = 1; $i -le $ws.Cells.FormatConditions.Count; $i++) {
    $fc = $ws.Cells.FormatConditions.Item($i)
    $fc.ModifyAppliesToRange($ws.Range("B2:H2"))
}

# Update the view: selection moves to C16, and the sheet is scrolled
# back so there's no special top-left cell (matches the post-edit state).
$ws.Range("A1").Select()
$ws.Range("C16").Select()
